# eooError added into core and updated
# Add a new "new speed" column (E) to the speed-test sheet:
#  - E1 header label "new speed"
#  - E3 value 9
# and leave the selection on the newly added cell, matching the author's
# final cursor position (E4) when they saved the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E1").Value = "new speed"
$ws.Range("E3").Value = 9

$ws.Activate()
$ws.Range("E4").Select()
